$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.431.29'
$ws.Range('E2').Value = '  -0.14%  '
$ws.Range('D3').Value = '1.726.51'
$ws.Range('E3').Value = '  +0.06%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9992'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '243.89'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.31%  '
$ws.Range('E6').Value = '  -0.06%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4913'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +1.84%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2616'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -1.84%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06203'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +0.39%  '
$ws.Range('D10').Value = '1.726.68'
$ws.Range('E10').Value = '  +0.02%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07022'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -2.29%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.54'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +0.01%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.567'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +1.12%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6019'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -1.22%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '77.34'
$ws.Range('D15').ClearFormats()
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.9996'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -0.05%  '
$ws.Range('D17').Value = '26.435.22'
$ws.Range('E17').Value = '  -0.17%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.9993'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -0.08%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007222'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +3.99%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.37'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.90%  '
$ws.Range('D21').Value = '1.944.57'
$ws.Range('E21').Value = '  -0.50%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.477'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.72%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.599'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -1.88%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.172'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -1.28%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '137.62'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.35%  '
$ws.Range('E26').Value = '  -0.48%  '
$ws.Range('B27').Value = 'BitcoinCash'
$ws.Range('C27').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '107.14'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.59%  '
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.387'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -0.69%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.706'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -3.95%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.965'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +0.18%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.07982'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.34%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.678'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -0.07%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04531'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.59%  '
$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.601'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.56%  '
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9999'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +0.89%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6279'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +0.50%  '
$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9128'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +0.03%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.968'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -4.44%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.390'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +0.60%  '
$ws.Range('B40').Value = 'PaxDollar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.000'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -0.52%  '
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.01486'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -1.00%  '
$ws.Range('B42').Value = 'Quant'
$ws.Range('C42').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '99.94'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -3.58%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.442'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -3.46%  '
$ws.Range('B44').Value = 'TheSandbox'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.3856'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.04%  '
$ws.Range('B45').Value = 'Aptos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '6.726'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -2.77%  '
$ws.Range('B46').Value = 'Algorand'
$ws.Range('C46').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.1155'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -1.87%  '
$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05366'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.39%  '
$ws.Range('B48').Value = 'Elrond'
$ws.Range('C48').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '30.14'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -0.77%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.701'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.70%  '
$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.240'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.64%  '
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '51.04'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +0.08%  '
